# Foundation Filed Cases, Related Cases, Not Litigated
#
# "Related Cases" sheet: the data row describing the Ori Spado / RICO Act
# case (previously the last row of the first case block, row 24) moves up
# to become the new row 5. Rows 5-23 shift down to become rows 6-24; the
# rows below (25-27) are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Related Cases")

# 1) Insert a fresh blank row at position 5; rows 5-23 shift down to 6-24,
#    and the row that used to be 24 is now (temporarily) duplicated at 25.
$ws.Rows.Item(5).Insert()

# 2) Populate the new row 5 with the data that used to live in row 24.
$ws.Range("A5").Value = "Criminal, Implications"
$ws.Range("B5").Value = "Artists, Witnesses, Hollywood"
$ws.Range("C5").Value = "2000s-2010"
$ws.Range("D5").Value = "Artists, Witnesses, Nightclubs, Singers"
$ws.Range("E5").Value = "California"
$ws.Range("F5").Value = "United States"
$ws.Range("G5").Value = "Ori Spado"
$ws.Range("J5").Value = "RICO Act, Association Crime Enterprises"
$ws.Range("K5").Value = "Informant, Hollywood"
$ws.Range("M5").Value = "FBI, Federal Government"
$ws.Range("P5").Value = "Implications, Implicated Hollywood Fixing, Implicating Crimes"
$ws.Range("Q5").Value = "Ori Spado, Witnesses, Book, Interviews"
$ws.Range("S5").Value = "https://theaccidentalgangster.com/"
$ws.Range("T5").Value = "https://www.amazon.com/ACCIDENTAL-GANGSTER-Insurance-Salesman-Hollywood/dp/1948239469"

# 3) The moved row keeps its original custom row height.
$ws.Rows.Item(5).RowHeight = 15.75

# 4) Remove the now-duplicated old row (pushed down to 25 by the insert) -
#    its content now lives in row 5.
$ws.Rows.Item(25).Delete()

# 5) The four existing hyperlinks are anchored to cells that just shifted
#    down one row (S22->S23, S7->S8, T13->T14, T17->T18). This runtime
#    doesn't re-anchor <hyperlinks> entries when rows are inserted/deleted,
#    and per-item Hyperlink.Delete()/.Range= are no-ops for links loaded
#    from the file, so clear the whole collection and re-add all four at
#    their correct post-shift addresses, keeping the original URLs/order.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("S23"), "https://people.com/movies/kip-pardue-fined-sexual-misconduct/")
$ws.Hyperlinks.Add($ws.Range("S8"), "https://en.wikipedia.org/wiki/Marc_Collins-Rector")
$ws.Hyperlinks.Add($ws.Range("T14"), "https://people.com/tv/corey-haim-mother-names-his-alleged-abuser/")
$ws.Hyperlinks.Add($ws.Range("T18"), "https://www.rollingstone.com/culture/culture-news/bryan-singer-abuse-allegations-blake-stuerman-1274239/")
